$wb = $excel.ActiveWorkbook

# "Test Results" sheet is the 4th sheet (Device, Zone, Category, Test Results)
$ws = $wb.Worksheets.Item("Test Results")

# Rows 16-24 correspond to Device test results (D01-D09) for
# Create / Read / Update / Delete tests. Mark them all as passed (TRUE),
# reflecting that test cases were added and results were written.
$ws.Range("B16:E24").Value = $true

# Select the full used range on the sheet, matching the resulting selection
# state stored in the workbook.
$ws.Activate()
$ws.Range("A1:E24").Select()
